# Generate Report for Handback
# Updates the localization-status workbook so that the "zh-cn" and "de-de"
# sheets show the markdown source file + generated XLIFF handback file in
# the "Latest Target File" / "Latest Handback File" columns, stamps the
# "Latest Handback DateTime" column, and flips the Status column from
# "Ready for handoff" to "Handed back: in sync with en-US" everywhere that
# text is used (Overview + both language sheets share the string).

$wb = $excel.ActiveWorkbook

$mdFile      = "7b1e7c19-20d0-4004-8d1a-136c09bbf563.md"
$mdUrl       = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/17ddc9232c9faea7f520ca5cf4420882fb868f77/e2e/7b1e7c19-20d0-4004-8d1a-136c09bbf563.md"
$zhXlf       = "7b1e7c19-20d0-4004-8d1a-136c09bbf563.7bd3990a85e2cc549418b1a94481ffd36c08d550.zh-cn.xlf"
$deXlf       = "7b1e7c19-20d0-4004-8d1a-136c09bbf563.7bd3990a85e2cc549418b1a94481ffd36c08d550.de-de.xlf"
$zhHandback  = "2016-08-30 15:22:26"
$deHandback  = "2016-08-30 15:22:33"
$statusText  = "Handed back: in sync with en-US"

# The status text used to be "Ready for handoff"; it now reads
# "Handed back: in sync with en-US" on every sheet that references it
# (Overview!E2:F3 and the Status column of the two language sheets).
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value2 -eq "Ready for handoff") {
                $cell.Value = $statusText
            }
        }
    }
}

foreach ($info in @(
        @{ Sheet = "zh-cn"; Xlf = $zhXlf; Handback = $zhHandback },
        @{ Sheet = "de-de"; Xlf = $deXlf; Handback = $deHandback }
    )) {

    $ws = $wb.Worksheets.Item($info.Sheet)

    foreach ($row in 2, 3) {
        # Column I: "Latest Target File" -> link to the markdown source file.
        $iCell = $ws.Cells.Item($row, 9)
        $iCell.Value = $mdFile
        $ws.Hyperlinks.Add($iCell, $mdUrl, "", "", $mdFile)
        $iCell.Style = "HyperLink"

        # Column J: "Latest Handback File" -> generated XLIFF file name.
        $ws.Cells.Item($row, 10).Value = $info.Xlf

        # Column K: "Latest Handback DateTime" -> timestamp of the handback.
        $ws.Cells.Item($row, 11).Value = $info.Handback
    }

    # Columns that now hold longer strings get wider to fit the content,
    # matching the report generator's auto-sizing behaviour.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

# Overview sheet's zh-cn/de-de status columns also grow to fit the new text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

Write-Output "Generate Report for Handback: done"
